$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 6.305023228865777
$ws.Range("R2").Value = 56.745209059792
$ws.Range("S2").Value = 0.01950610157132724
$ws.Range("T2").Value = 0.01950610157132724

$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("S3").Value = 0.4425091540123311
$ws.Range("T3").Value = 0.4425091540123311

$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("R4").Value = 704.8962577378039
$ws.Range("S4").Value = 0.2423072930473132
$ws.Range("T4").Value = 0.2423072930473132

$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 2.646874223688
$ws.Range("R5").Value = 23.821868013192
$ws.Range("S5").Value = 0.008188740244034582
$ws.Range("T5").Value = 0.008188740244034582

$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("S6").Value = 0.1857671305854846
$ws.Range("T6").Value = 0.1857671305854845

$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("S7").Value = 0.1017215805395094
$ws.Range("T7").Value = 0.1017215805395094
